$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 761.3333
$ws.Range("J19").Value = 442.85715
$ws.Range("L19").Value = 442.85715
$ws.Range("N19").Value = -792.85715
$ws.Range("H69").Value = 5866.7896
$ws.Range("J69").Value = 6081.6665
$ws.Range("L69").Value = 18244.9995
$ws.Range("N69").Value = -19992.9995
$ws.Range("H72").Value = 5866.7896
$ws.Range("J72").Value = 6081.6665
$ws.Range("L72").Value = 54734.9985
$ws.Range("N72").Value = -63470.9985
$ws.Range("H76").Value = 3632.8333
$ws.Range("I76").Value = 3999.3333
$ws.Range("J76").Value = 3266.3333
$ws.Range("K76").Value = 3999.3333
$ws.Range("L76").Value = 3266.3333
$ws.Range("M76").Value = -3684.3333
$ws.Range("N76").Value = -3896.3333
$ws.Range("H79").Value = 3632.8333
$ws.Range("I79").Value = 3999.3333
$ws.Range("J79").Value = 3266.3333
$ws.Range("K79").Value = 3999.3333
$ws.Range("L79").Value = 3266.3333
$ws.Range("M79").Value = -2907.3333
$ws.Range("N79").Value = -5450.3333
$ws.Range("H92").Value = 106.818184
$ws.Range("I92").Value = 119
$ws.Range("J92").Value = 52
$ws.Range("K92").Value = 119
$ws.Range("L92").Value = 52
$ws.Range("M92").Value = 1129
$ws.Range("N92").Value = -2548
$ws.Range("H113").Value = 2863.5
$ws.Range("I113").Value = 484.66666
$ws.Range("K113").Value = 484.66666
$ws.Range("M113").Value = 2769.33334
$ws.Range("H132").Value = 9536.786
$ws.Range("J132").Value = 11301.25
$ws.Range("L132").Value = 33903.75
$ws.Range("N132").Value = -38963.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H41").Value = 138
$ws.Range("I41").Value = 147.75
$ws.Range("J41").Value = 99
$ws.Range("K41").Value = 147.75
$ws.Range("L41").Value = 99
$ws.Range("M41").Value = 266.25
$ws.Range("N41").Value = -927
$ws.Range("H45").Value = 2840.4443
$ws.Range("I45").Value = 1576
$ws.Range("J45").Value = 4421
$ws.Range("K45").Value = 1576
$ws.Range("L45").Value = 4421
$ws.Range("M45").Value = -1199
$ws.Range("N45").Value = -5175
$ws.Range("H63").Value = 13742
$ws.Range("I63").Value = 1898.3334
$ws.Range("J63").Value = 22624.75
$ws.Range("K63").Value = 1898.3334
$ws.Range("L63").Value = 22624.75
$ws.Range("M63").Value = -1212.3334
$ws.Range("N63").Value = -23996.75
$ws.Range("H66").Value = 13742
$ws.Range("I66").Value = 1898.3334
$ws.Range("J66").Value = 22624.75
$ws.Range("K66").Value = 9491.666999999999
$ws.Range("L66").Value = 113123.75
$ws.Range("M66").Value = -6059.666999999999
$ws.Range("N66").Value = -119987.75
$ws.Range("H132").Value = 3979.9285
$ws.Range("I132").Value = 3876.182
$ws.Range("K132").Value = 11628.546
$ws.Range("M132").Value = -9098.545999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1600.1111
$ws.Range("I22").Value = 1600.1111
$ws.Range("K22").Value = 1600.1111
$ws.Range("M22").Value = -1427.1111
$ws.Range("H49").Value = 50000
$ws.Range("J49").Value = 50000
$ws.Range("L49").Value = 50000
$ws.Range("N49").Value = -50478
$ws.Range("H86").Value = 5059.375
$ws.Range("I86").Value = 3930.8
$ws.Range("J86").Value = 6940.3335
$ws.Range("K86").Value = 3930.8
$ws.Range("L86").Value = 6940.3335
$ws.Range("M86").Value = -2807.8
$ws.Range("N86").Value = -9186.333500000001
$ws.Range("H89").Value = 5059.375
$ws.Range("I89").Value = 3930.8
$ws.Range("J89").Value = 6940.3335
$ws.Range("K89").Value = 19654
$ws.Range("L89").Value = 34701.6675
$ws.Range("M89").Value = -14038
$ws.Range("N89").Value = -45933.6675
$ws.Range("H105").Value = 1890
$ws.Range("I105").Value = 1686.6666
$ws.Range("K105").Value = 1686.6666
$ws.Range("M105").Value = 60.33339999999998
$ws.Range("H107").Value = 5214.067
$ws.Range("I107").Value = 4421.2
$ws.Range("K107").Value = 4421.2
$ws.Range("M107").Value = -2501.2
$ws.Range("H134").Value = 3399.7778
$ws.Range("I134").Value = 3399.7778
$ws.Range("K134").Value = 10199.3334
$ws.Range("M134").Value = -7664.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 449996.25
$ws.Range("J4").Value = 708000
$ws.Range("L4").Value = 708000
$ws.Range("N4").Value = -708224
$ws.Range("H7").Value = 493.125
$ws.Range("I7").Value = 95
$ws.Range("K7").Value = 95
$ws.Range("M7").Value = 18
$ws.Range("H21").Value = 10500
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 10500
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 10500
$ws.Range("M21").ClearContents()
$ws.Range("N21").Value = -10970
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H58").Value = 3455.5715
$ws.Range("I58").Value = 3029.4614
$ws.Range("K58").Value = 3029.4614
$ws.Range("M58").Value = -2826.4614
$ws.Range("H62").Value = 5012.2
$ws.Range("I62").Value = 4761.25
$ws.Range("J62").Value = 6016
$ws.Range("K62").Value = 4761.25
$ws.Range("L62").Value = 6016
$ws.Range("M62").Value = -4137.25
$ws.Range("N62").Value = -7264
$ws.Range("H65").Value = 5012.2
$ws.Range("I65").Value = 4761.25
$ws.Range("J65").Value = 6016
$ws.Range("K65").Value = 23806.25
$ws.Range("L65").Value = 30080
$ws.Range("M65").Value = -20686.25
$ws.Range("N65").Value = -36320
$ws.Range("H94").Value = 6914
$ws.Range("I94").Value = 4257
$ws.Range("J94").Value = 11342.333
$ws.Range("K94").Value = 4257
$ws.Range("L94").Value = 11342.333
$ws.Range("M94").Value = -3806
$ws.Range("N94").Value = -12244.333
$ws.Range("H136").Value = 3455.5715
$ws.Range("I136").Value = 3029.4614
$ws.Range("K136").Value = 9088.3842
$ws.Range("M136").Value = -6538.3842

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 97.45
$ws.Range("I12").Value = 12
$ws.Range("J12").Value = 125.933334
$ws.Range("K12").Value = 36
$ws.Range("L12").Value = 377.800002
$ws.Range("M12").Value = 137
$ws.Range("N12").Value = -723.8000019999999
$ws.Range("H26").Value = 43.6
$ws.Range("J26").Value = 40
$ws.Range("L26").Value = 120
$ws.Range("N26").Value = -696
$ws.Range("H52").Value = 1574
$ws.Range("J52").Value = 1574
$ws.Range("L52").Value = 4722
$ws.Range("N52").Value = -5254
$ws.Range("H81").Value = 1140
$ws.Range("J81").Value = 2000
$ws.Range("L81").Value = 6000
$ws.Range("N81").Value = -8246
$ws.Range("H84").Value = 1140
$ws.Range("J84").Value = 2000
$ws.Range("L84").Value = 18000
$ws.Range("N84").Value = -29232
$ws.Range("H86").Value = 410.33334
$ws.Range("I86").Value = 428
$ws.Range("J86").Value = 348.5
$ws.Range("K86").Value = 1284
$ws.Range("L86").Value = 1045.5
$ws.Range("M86").Value = -98
$ws.Range("N86").Value = -3417.5
$ws.Range("H89").Value = 410.33334
$ws.Range("I89").Value = 428
$ws.Range("J89").Value = 348.5
$ws.Range("K89").Value = 3852
$ws.Range("L89").Value = 3136.5
$ws.Range("M89").Value = 2076
$ws.Range("N89").Value = -14992.5
$ws.Range("H104").Value = 8902.308000000001
$ws.Range("J104").Value = 9180
$ws.Range("L104").Value = 27540
$ws.Range("N104").Value = -32782
$ws.Range("H113").Value = 908.4375
$ws.Range("I113").Value = 620.5
$ws.Range("J113").Value = 949.5714
$ws.Range("K113").Value = 1861.5
$ws.Range("L113").Value = 2848.7142
$ws.Range("M113").Value = 308.5
$ws.Range("N113").Value = -7188.7142

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 7131.846
$ws.Range("I113").Value = 6062
$ws.Range("K113").Value = 6062
$ws.Range("M113").Value = -3892
$ws.Range("H122").Value = 4444
$ws.Range("I122").Value = 2777.5
$ws.Range("K122").Value = 8332.5
$ws.Range("M122").Value = -5882.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 666.6667
$ws.Range("I22").Value = 666.6667
$ws.Range("K22").Value = 666.6667
$ws.Range("M22").Value = -371.6667
$ws.Range("H24").Value = 2000
$ws.Range("I24").Value = 2000
$ws.Range("K24").Value = 2000
$ws.Range("M24").Value = -1657
$ws.Range("H27").Value = 666.6667
$ws.Range("I27").Value = 666.6667
$ws.Range("K27").Value = 666.6667
$ws.Range("M27").Value = -559.6667
$ws.Range("H55").Value = 1016
$ws.Range("I55").Value = 1016
$ws.Range("K55").Value = 1016
$ws.Range("M55").Value = -843

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 4467.8
$ws.Range("I122").Value = 2320
$ws.Range("J122").Value = 5899.6665
$ws.Range("K122").Value = 6960
$ws.Range("L122").Value = 17698.9995
$ws.Range("M122").Value = -4510
$ws.Range("N122").Value = -22598.9995
$ws.Range("H126").Value = 5993.647
$ws.Range("I126").Value = 3737.75
$ws.Range("K126").Value = 11213.25
$ws.Range("M126").Value = -8743.25
$ws.Range("H129").Value = 239499.5
$ws.Range("I129").Value = 49999
$ws.Range("K129").Value = 49999
$ws.Range("M129").Value = -44999
